$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1562.625
$ws.Range("I33").Value = 1562.625
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 1562.625
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -1333.625
$ws.Range("H40").Value = 4311.875
$ws.Range("I40").Value = 3797.5
$ws.Range("K40").Value = 3797.5
$ws.Range("M40").Value = -3622.5
$ws.Range("H62").Value = 148133.28
$ws.Range("I62").Value = 203997.6
$ws.Range("K62").Value = 203997.6
$ws.Range("M62").Value = -203373.6
$ws.Range("H65").Value = 148133.28
$ws.Range("I65").Value = 203997.6
$ws.Range("K65").Value = 1019988
$ws.Range("M65").Value = -1016868
$ws.Range("H76").Value = 3842.1538
$ws.Range("I76").Value = 3329
$ws.Range("K76").Value = 3329
$ws.Range("M76").Value = -3014
$ws.Range("H79").Value = 3842.1538
$ws.Range("I79").Value = 3329
$ws.Range("K79").Value = 3329
$ws.Range("M79").Value = -2237
$ws.Range("H125").Value = 33643.8
$ws.Range("I125").Value = 107999
$ws.Range("J125").Value = 1777.2858
$ws.Range("K125").Value = 971991
$ws.Range("L125").Value = 15995.5722
$ws.Range("M125").Value = -969531
$ws.Range("N125").Value = -20915.5722
$ws.Range("H132").Value = 4642.9287
$ws.Range("I132").Value = 2037.8889
$ws.Range("K132").Value = 6113.6667
$ws.Range("M132").Value = -3583.6667
$ws.Range("H137").Value = 3447.42
$ws.Range("I137").Value = 1400.4147
$ws.Range("J137").Value = 12772.667
$ws.Range("K137").Value = 4201.2441
$ws.Range("L137").Value = 38318.001
$ws.Range("M137").Value = -1651.2441
$ws.Range("N137").Value = -43418.001
$ws.Range("H138").Value = 2464.1428
$ws.Range("I138").Value = 1715.8667
$ws.Range("J138").Value = 3327.5386
$ws.Range("K138").Value = 5147.6001
$ws.Range("L138").Value = 9982.6158
$ws.Range("M138").Value = -7.600099999999657
$ws.Range("N138").Value = -20262.6158
$ws.Range("H141").Value = 48300
$ws.Range("J141").Value = 5250
$ws.Range("L141").Value = 15750
$ws.Range("N141").Value = -26110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1720.878
$ws.Range("I2").Value = 1461.6774
$ws.Range("J2").Value = 2524.4
$ws.Range("K2").Value = 1461.6774
$ws.Range("L2").Value = 2524.4
$ws.Range("M2").Value = -1348.6774
$ws.Range("N2").Value = -2750.4
$ws.Range("H45").Value = 462017.22
$ws.Range("I45").Value = 844631.75
$ws.Range("J45").Value = 2879.8
$ws.Range("K45").Value = 844631.75
$ws.Range("L45").Value = 2879.8
$ws.Range("M45").Value = -844254.75
$ws.Range("N45").Value = -3633.8
$ws.Range("H61").Value = 1165.0571
$ws.Range("I61").Value = 783.75
$ws.Range("K61").Value = 783.75
$ws.Range("M61").Value = -571.75
$ws.Range("H116").Value = 1720.878
$ws.Range("I116").Value = 1461.6774
$ws.Range("J116").Value = 2524.4
$ws.Range("K116").Value = 1461.6774
$ws.Range("L116").Value = 2524.4
$ws.Range("M116").Value = 832.3226
$ws.Range("N116").Value = -7112.4
$ws.Range("H122").Value = 2333.1538
$ws.Range("I122").Value = 2195.5557
$ws.Range("K122").Value = 6586.6671
$ws.Range("M122").Value = -4136.6671
$ws.Range("H136").Value = 1165.0571
$ws.Range("I136").Value = 783.75
$ws.Range("K136").Value = 2351.25
$ws.Range("M136").Value = 198.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1720.878
$ws.Range("I3").Value = 1461.6774
$ws.Range("J3").Value = 2524.4
$ws.Range("K3").Value = 1461.6774
$ws.Range("L3").Value = 2524.4
$ws.Range("M3").Value = -1347.6774
$ws.Range("N3").Value = -2752.4
$ws.Range("H80").Value = 1171.625
$ws.Range("I80").Value = 735.8333
$ws.Range("J80").Value = 1433.1
$ws.Range("K80").Value = 735.8333
$ws.Range("L80").Value = 1433.1
$ws.Range("M80").Value = 262.1667
$ws.Range("N80").Value = -3429.1
$ws.Range("H83").Value = 1171.625
$ws.Range("I83").Value = 735.8333
$ws.Range("J83").Value = 1433.1
$ws.Range("K83").Value = 3679.1665
$ws.Range("L83").Value = 7165.5
$ws.Range("M83").Value = 1312.8335
$ws.Range("N83").Value = -17149.5
$ws.Range("H99").Value = 1930.0769
$ws.Range("I99").Value = 1372.8182
$ws.Range("J99").Value = 4995
$ws.Range("K99").Value = 1372.8182
$ws.Range("L99").Value = 4995
$ws.Range("M99").Value = 125.1818000000001
$ws.Range("N99").Value = -7991
$ws.Range("H105").Value = 9402.700000000001
$ws.Range("I105").Value = 10461.143
$ws.Range("K105").Value = 10461.143
$ws.Range("M105").Value = -8714.143
$ws.Range("H134").Value = 1699.9445
$ws.Range("I134").Value = 1614.4286
$ws.Range("J134").Value = 1999.25
$ws.Range("K134").Value = 4843.2858
$ws.Range("L134").Value = 5997.75
$ws.Range("M134").Value = -2308.2858
$ws.Range("N134").Value = -11067.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 220.45454
$ws.Range("I22").Value = 215
$ws.Range("K22").Value = 215
$ws.Range("M22").Value = 135
$ws.Range("H31").Value = 2100.9656
$ws.Range("I31").Value = 1551.591
$ws.Range("J31").Value = 3827.5715
$ws.Range("K31").Value = 1551.591
$ws.Range("L31").Value = 3827.5715
$ws.Range("M31").Value = -1256.591
$ws.Range("N31").Value = -4417.5715
$ws.Range("H34").Value = 2100.9656
$ws.Range("I34").Value = 1551.591
$ws.Range("J34").Value = 3827.5715
$ws.Range("K34").Value = 1551.591
$ws.Range("L34").Value = 3827.5715
$ws.Range("M34").Value = -1349.591
$ws.Range("N34").Value = -4231.5715
$ws.Range("H87").Value = 25330
$ws.Range("J87").Value = 25330
$ws.Range("L87").Value = 25330
$ws.Range("N87").Value = -27702
$ws.Range("H90").Value = 25330
$ws.Range("J90").Value = 25330
$ws.Range("L90").Value = 75990
$ws.Range("N90").Value = -87846
$ws.Range("H132").Value = 2542.5715
$ws.Range("I132").Value = 2381
$ws.Range("K132").Value = 7143
$ws.Range("M132").Value = -4613
$ws.Range("H134").Value = 4496.5386
$ws.Range("I134").Value = 4454.5835
$ws.Range("K134").Value = 13363.7505
$ws.Range("M134").Value = -10828.7505

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 111156500
$ws.Range("J37").Value = 111156500
$ws.Range("L37").Value = 333469500
$ws.Range("N37").Value = -333469724
$ws.Range("H39").Value = 3892.4285
$ws.Range("J39").Value = 4750
$ws.Range("L39").Value = 14250
$ws.Range("N39").Value = -14838
$ws.Range("H55").Value = 7814289.5
$ws.Range("J55").Value = 12502590
$ws.Range("L55").Value = 37507770
$ws.Range("N55").Value = -37508124
$ws.Range("H68").Value = 3248.75
$ws.Range("I68").Value = 1497.5
$ws.Range("K68").Value = 4492.5
$ws.Range("M68").Value = -3681.5
$ws.Range("H71").Value = 3248.75
$ws.Range("I71").Value = 1497.5
$ws.Range("K71").Value = 13477.5
$ws.Range("M71").Value = -9421.5
$ws.Range("H92").Value = 4999.5
$ws.Range("J92").Value = 4999.5
$ws.Range("L92").Value = 14998.5
$ws.Range("N92").Value = -17494.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4443.5
$ws.Range("I70").Value = 4443.5
$ws.Range("K70").Value = 4443.5
$ws.Range("M70").Value = -4173.5
$ws.Range("H73").Value = 4443.5
$ws.Range("I73").Value = 4443.5
$ws.Range("K73").Value = 4443.5
$ws.Range("M73").Value = -3507.5
$ws.Range("H122").Value = 3311.8823
$ws.Range("I122").Value = 3393.875
$ws.Range("K122").Value = 10181.625
$ws.Range("M122").Value = -7731.625
$ws.Range("H132").Value = 1951
$ws.Range("I132").Value = 1812.25
$ws.Range("K132").Value = 5436.75
$ws.Range("M132").Value = -2906.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 686131.4399999999
$ws.Range("I132").Value = 1003466.06
$ws.Range("J132").Value = 6128.7144
$ws.Range("K132").Value = 3010398.18
$ws.Range("L132").Value = 18386.1432
$ws.Range("M132").Value = -3007868.18
$ws.Range("N132").Value = -23446.1432
$ws.Range("H136").Value = 2942.913
$ws.Range("I136").Value = 2388.4443
$ws.Range("K136").Value = 7165.3329
$ws.Range("M136").Value = -4615.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5489.2
$ws.Range("I81").Value = 5489.2
$ws.Range("K81").Value = 10978.4
$ws.Range("M81").Value = -9917.4
$ws.Range("H84").Value = 5489.2
$ws.Range("I84").Value = 5489.2
$ws.Range("K84").Value = 54892
$ws.Range("M84").Value = -49588
$ws.Range("H96").Value = 7250
$ws.Range("J96").Value = 3500
$ws.Range("L96").Value = 3500
$ws.Range("N96").Value = -6246
$ws.Range("H132").Value = 40873.273
$ws.Range("I132").Value = 37271.42
$ws.Range("K132").Value = 111814.26
$ws.Range("M132").Value = -109284.26
